$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf13"
$ws.Range("C2").Value = "Tnfrsf14"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.637903
$ws.Range("H2").Value = 1.913709
$ws.Range("I2").Value = 0.1229013127714845
$ws.Range("J2").Value = 0.1229013127714844
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.667667333333333
$ws.Range("N2").Value = 11.003002
$ws.Range("O2").Value = 0.1424137080579054
$ws.Range("P2").Value = 0.1424137080579054
$ws.Range("Q2").Value = 2.339615994935333
$ws.Range("R2").Value = 21.056543954418
$ws.Range("S2").Value = 0.01750283167697151
$ws.Range("T2").Value = 0.0175028316769715

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfsf13"
$ws.Range("C3").Value = "Tnfrsf14"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.637903
$ws.Range("H3").Value = 1.913709
$ws.Range("I3").Value = 0.1229013127714845
$ws.Range("J3").Value = 0.1229013127714844
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.497702
$ws.Range("N3").Value = 19.493106
$ws.Range("O3").Value = 0.2523025540689536
$ws.Range("P3").Value = 0.2523025540689536
$ws.Range("Q3").Value = 4.144903598906001
$ws.Range("R3").Value = 37.304132390154
$ws.Range("S3").Value = 0.03100831511067283
$ws.Range("T3").Value = 0.03100831511067283

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tnfsf13"
$ws.Range("C4").Value = "Tnfrsf14"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.637903
$ws.Range("H4").Value = 1.913709
$ws.Range("I4").Value = 0.1229013127714845
$ws.Range("J4").Value = 0.1229013127714844
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.399313
$ws.Range("N4").Value = 43.197939
$ws.Range("O4").Value = 0.5591182000556945
$ws.Range("P4").Value = 0.5591182000556945
$ws.Range("Q4").Value = 9.185364960639
$ws.Range("R4").Value = 82.66828464575099
$ws.Range("S4").Value = 0.06871636078127433
$ws.Range("T4").Value = 0.06871636078127431

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tnfsf13"
$ws.Range("C5").Value = "Tnfrsf14"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.637903
$ws.Range("H5").Value = 1.913709
$ws.Range("I5").Value = 0.1229013127714845
$ws.Range("J5").Value = 0.1229013127714844
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.188929333333334
$ws.Range("N5").Value = 3.566788
$ws.Range("O5").Value = 0.04616553781744659
$ws.Range("P5").Value = 0.04616553781744658
$ws.Range("Q5").Value = 0.7584215885213335
$ws.Range("R5").Value = 6.825794296692
$ws.Range("S5").Value = 0.005673805202565797
$ws.Range("T5").Value = 0.005673805202565795

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tnfsf13"
$ws.Range("C6").Value = "Tnfrsf14"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.170281
$ws.Range("H6").Value = 0.510843
$ws.Range("I6").Value = 0.03280711713229307
$ws.Range("J6").Value = 0.03280711713229307
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.667667333333333
$ws.Range("N6").Value = 11.003002
$ws.Range("O6").Value = 0.1424137080579054
$ws.Range("P6").Value = 0.1424137080579054
$ws.Range("Q6").Value = 0.6245340611873332
$ws.Range("R6").Value = 5.620806550686
$ws.Range("S6").Value = 0.004672183201499892
$ws.Range("T6").Value = 0.004672183201499891

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tnfsf13"
$ws.Range("C7").Value = "Tnfrsf14"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.170281
$ws.Range("H7").Value = 0.510843
$ws.Range("I7").Value = 0.03280711713229307
$ws.Range("J7").Value = 0.03280711713229307
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.497702
$ws.Range("N7").Value = 19.493106
$ws.Range("O7").Value = 0.2523025540689536
$ws.Range("P7").Value = 0.2523025540689536
$ws.Range("Q7").Value = 1.106435194262
$ws.Range("R7").Value = 9.957916748358002
$ws.Range("S7").Value = 0.008277319444116865
$ws.Range("T7").Value = 0.008277319444116865

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Tnfsf13"
$ws.Range("C8").Value = "Tnfrsf14"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.170281
$ws.Range("H8").Value = 0.510843
$ws.Range("I8").Value = 0.03280711713229307
$ws.Range("J8").Value = 0.03280711713229307
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.399313
$ws.Range("N8").Value = 43.197939
$ws.Range("O8").Value = 0.5591182000556945
$ws.Range("P8").Value = 0.5591182000556945
$ws.Range("Q8").Value = 2.451929416953
$ws.Range("R8").Value = 22.067364752577
$ws.Range("S8").Value = 0.01834305628002404
$ws.Range("T8").Value = 0.01834305628002404

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Tnfsf13"
$ws.Range("C9").Value = "Tnfrsf14"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.170281
$ws.Range("H9").Value = 0.510843
$ws.Range("I9").Value = 0.03280711713229307
$ws.Range("J9").Value = 0.03280711713229307
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.188929333333334
$ws.Range("N9").Value = 3.566788
$ws.Range("O9").Value = 0.04616553781744659
$ws.Range("P9").Value = 0.04616553781744658
$ws.Range("Q9").Value = 0.2024520758093334
$ws.Range("R9").Value = 1.822068682284
$ws.Range("S9").Value = 0.001514558206652276
$ws.Range("T9").Value = 0.001514558206652275

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Tnfsf13"
$ws.Range("C10").Value = "Tnfrsf14"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.407124666666667
$ws.Range("H10").Value = 10.221374
$ws.Range("I10").Value = 0.6564322386153376
$ws.Range("J10").Value = 0.6564322386153377
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.667667333333333
$ws.Range("N10").Value = 11.003002
$ws.Range("O10").Value = 0.1424137080579054
$ws.Range("P10").Value = 0.1424137080579054
$ws.Range("Q10").Value = 12.49619984052755
$ws.Range("R10").Value = 112.465798564748
$ws.Range("S10").Value = 0.09348494918996199
$ws.Range("T10").Value = 0.09348494918996199

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Tnfsf13"
$ws.Range("C11").Value = "Tnfrsf14"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.407124666666667
$ws.Range("H11").Value = 10.221374
$ws.Range("I11").Value = 0.6564322386153376
$ws.Range("J11").Value = 0.6564322386153377
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.497702
$ws.Range("N11").Value = 19.493106
$ws.Range("O11").Value = 0.2523025540689536
$ws.Range("P11").Value = 0.2523025540689536
$ws.Range("Q11").Value = 22.13848076084934
$ws.Range("R11").Value = 199.246326847644
$ws.Range("S11").Value = 0.1656195303758504
$ws.Range("T11").Value = 0.1656195303758505

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Tnfsf13"
$ws.Range("C12").Value = "Tnfrsf14"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.407124666666667
$ws.Range("H12").Value = 10.221374
$ws.Range("I12").Value = 0.6564322386153376
$ws.Range("J12").Value = 0.6564322386153377
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 14.399313
$ws.Range("N12").Value = 43.197939
$ws.Range("O12").Value = 0.5591182000556945
$ws.Range("P12").Value = 0.5591182000556945
$ws.Range("Q12").Value = 49.060254505354
$ws.Range("R12").Value = 441.542290548186
$ws.Range("S12").Value = 0.3670232117131377
$ws.Range("T12").Value = 0.3670232117131378

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Tnfsf13"
$ws.Range("C13").Value = "Tnfrsf14"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.407124666666667
$ws.Range("H13").Value = 10.221374
$ws.Range("I13").Value = 0.6564322386153376
$ws.Range("J13").Value = 0.6564322386153377
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188929333333334
$ws.Range("N13").Value = 3.566788
$ws.Range("O13").Value = 0.04616553781744659
$ws.Range("P13").Value = 0.04616553781744658
$ws.Range("Q13").Value = 4.050830458523556
$ws.Range("R13").Value = 36.45747412671201
$ws.Range("S13").Value = 0.03030454733638749
$ws.Range("T13").Value = 0.03030454733638749

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Tnfsf13"
$ws.Range("C14").Value = "Tnfrsf14"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.975059
$ws.Range("H14").Value = 2.925177
$ws.Range("I14").Value = 0.1878593314808848
$ws.Range("J14").Value = 0.1878593314808848
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.667667333333333
$ws.Range("N14").Value = 11.003002
$ws.Range("O14").Value = 0.1424137080579054
$ws.Range("P14").Value = 0.1424137080579054
$ws.Range("Q14").Value = 3.576192042372666
$ws.Range("R14").Value = 32.185728381354
$ws.Range("S14").Value = 0.026753743989472
$ws.Range("T14").Value = 0.026753743989472

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Tnfsf13"
$ws.Range("C15").Value = "Tnfrsf14"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.975059
$ws.Range("H15").Value = 2.925177
$ws.Range("I15").Value = 0.1878593314808848
$ws.Range("J15").Value = 0.1878593314808848
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 6.497702
$ws.Range("N15").Value = 19.493106
$ws.Range("O15").Value = 0.2523025540689536
$ws.Range("P15").Value = 0.2523025540689536
$ws.Range("Q15").Value = 6.335642814418001
$ws.Range("R15").Value = 57.02078532976201
$ws.Range("S15").Value = 0.04739738913831341
$ws.Range("T15").Value = 0.04739738913831341

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Tnfsf13"
$ws.Range("C16").Value = "Tnfrsf14"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.975059
$ws.Range("H16").Value = 2.925177
$ws.Range("I16").Value = 0.1878593314808848
$ws.Range("J16").Value = 0.1878593314808848
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 14.399313
$ws.Range("N16").Value = 43.197939
$ws.Range("O16").Value = 0.5591182000556945
$ws.Range("P16").Value = 0.5591182000556945
$ws.Range("Q16").Value = 14.040179734467
$ws.Range("R16").Value = 126.361617610203
$ws.Range("S16").Value = 0.1050355712812584
$ws.Range("T16").Value = 0.1050355712812584

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Tnfsf13"
$ws.Range("C17").Value = "Tnfrsf14"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.975059
$ws.Range("H17").Value = 2.925177
$ws.Range("I17").Value = 0.1878593314808848
$ws.Range("J17").Value = 0.1878593314808848
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.188929333333334
$ws.Range("N17").Value = 3.566788
$ws.Range("O17").Value = 0.04616553781744659
$ws.Range("P17").Value = 0.04616553781744658
$ws.Range("Q17").Value = 1.159276246830667
$ws.Range("R17").Value = 10.433486221476
$ws.Range("S17").Value = 0.008672627071841021
$ws.Range("T17").Value = 0.008672627071841021

